$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two rows (old rows 11 and 12); remaining rows 7-10 will be
# overwritten below with the "rolled forward" data set.
$ws.Rows("11:12").Delete()

# Columns A (Caso) and B (F. De Reclamo) store plain text in this sheet, not
# numbers/dates, so force text formatting before assigning to keep them as
# literal strings (e.g. "6557", "8/4/2025") instead of being auto-converted.
$ws.Range("A7:B10").NumberFormat = "@"

# Row 7: Caso 6557 - ALBERDI, JUAN BAUTISTA AV. 1091
$ws.Range("A7").Value = "6557"
$ws.Range("B7").Value = "8/4/2025"
$ws.Range("C7").Value = "ALBERDI, JUAN BAUTISTA AV. 1091"
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = "Pedido a ADM"
$ws.Range("H7").Value = "Cables a baja altura"
$ws.Range("J7").Value = '{"direccionesNormalizadas": [{"altura": 1091, "cod_calle": 1033, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.443538", "y": "-34.622890"}, "direccion": "ALBERDI, JUAN BAUTISTA AV. 1091, CABA", "nombre_calle": "ALBERDI, JUAN BAUTISTA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K7").Value = -58.443538
$ws.Range("L7").Value = -34.62289
$ws.Range("M7").Value = "Boedo"
$ws.Range("N7").Value = "Capital Sur"

# Row 8: Caso 6193 - POLA 591
$ws.Range("A8").Value = "6193"
$ws.Range("B8").Value = "8/4/2025"
$ws.Range("C8").Value = "POLA 591"
$ws.Range("D8").Value = 9
$ws.Range("E8").Value = "Pedido a ADM"
$ws.Range("H8").Value = "Cable colgando y enrollado en arbol"
$ws.Range("J8").Value = '{"direccionesNormalizadas": [{"altura": 591, "cod_calle": 17105, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.507385", "y": "-34.644479"}, "direccion": "POLA 591, CABA", "nombre_calle": "POLA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K8").Value = -58.507385
$ws.Range("L8").Value = -34.644479
$ws.Range("M8").Value = "Devoto"
$ws.Range("N8").Value = "Capital Norte"

# Row 9: Caso 6345 - ARTIGAS, JOSE GERVASIO, GRAL. 924
$ws.Range("A9").Value = "6345"
$ws.Range("B9").Value = "8/5/2025"
$ws.Range("C9").Value = "ARTIGAS, JOSE GERVASIO, GRAL. 924"
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = "Pedido a ADM"
$ws.Range("H9").Value = "Tendido a baja altura"
$ws.Range("J9").Value = '{"direccionesNormalizadas": [{"altura": 924, "cod_calle": 1125, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.468116", "y": "-34.620216"}, "direccion": "ARTIGAS, JOSE GERVASIO, GRAL. 924, CABA", "nombre_calle": "ARTIGAS, JOSE GERVASIO, GRAL.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K9").Value = -58.468116
$ws.Range("L9").Value = -34.620216
$ws.Range("M9").Value = "Boedo"
$ws.Range("N9").Value = "Capital Sur"

# Row 10: Caso 6568 - Carlos E. Pellegrini 6030
$ws.Range("A10").Value = "6568"
$ws.Range("B10").Value = "8/5/2025"
$ws.Range("C10").Value = "Carlos E. Pellegrini 6030"
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = "Pedido a ADM"
$ws.Range("H10").Value = "Tendido a muy baja altura se solicita retiro o levantarlo"
$ws.Range("J10").Value = '{"direccionesNormalizadas": [{"altura": 6030, "cod_calle": 17053, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.504789", "y": "-34.564505"}, "direccion": "PELLEGRINI, CARLOS E. 6030, CABA", "nombre_calle": "PELLEGRINI, CARLOS E.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K10").Value = -58.504789
$ws.Range("L10").Value = -34.564505
$ws.Range("M10").Value = "Paternal"
$ws.Range("N10").Value = "Capital Norte"
